# Converts a "RRGGBB" hex colour string into the packed BGR integer that
# PowerPoint's COM RGB properties expect (VBA RGB() macro: R + G*256 + B*65536).
function ConvertTo-VbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 5 switches from the custom pink "Table_0" table
#        style to the built-in "No Style, Table Grid" style. ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{26F89C1C-5F9E-4C52-92D2-CBEA04F409FD}")

# --- 2. The deck's theme colour scheme changes from the "Red Violet"
#        (Integral) palette to the standard "Office" palette. ---
$officeColours = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hyperlink
    "954F72"  # followed hyperlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-VbaRgb($officeColours[$i - 1])
}
